# Enhance csv module error handling / simplify code structure.
# Appends one new data row (row 55) to each of the four worksheets,
# mirroring the layout already used by the existing rows (e.g. row 54).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "MID_LFT_#1"
        A = 45841.46162037037
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x68"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 360
        I = 7
    },
    @{
        Sheet = "MID_LFT_#2"
        A = 45841.46162037037
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x60"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 352
        I = 25
    },
    @{
        Sheet = "MID_PLT_#1"
        A = 45841.46162037037
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x69"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 105
        I = 15
    },
    @{
        Sheet = "MID_PLT_#2"
        A = 45841.46162037037
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7E"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 126
        I = 9
    }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = $entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}

Write-Host "Appended row 55 to all four worksheets."
